$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) After the "Table" (Heading2) paragraph, insert three new bulleted
#    ListParagraph items (all sharing the same new numbered list).
# ---------------------------------------------------------------------------
$tableHeading = $d.Content.Find.Execute("Table", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tableRange = $d.Content
$tableRange.Find.Execute("Table", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tablePara = $tableRange.Paragraphs(1)

$tablePara.Range.InsertParagraphAfter()
$item1 = $tablePara.Next()
$item1.Style = "ListParagraph"
$item1.Range.Text = "List of items is retrieved from database and put in Vuex store when list is created. If we route to another page and come back, the list should not refresh. "

$gallery = $word.ListGalleries.Item(1)
$bulletTemplate = $gallery.ListTemplates.Item(1)
$item1.Range.ListFormat.ApplyListTemplate($bulletTemplate)

$item1.Range.InsertParagraphAfter()
$item2 = $item1.Next()
$item2.Style = "ListParagraph"
$item2.Range.Text = "List of items is maintained in store and is automatically refreshed every 10 minutes. "
$item2.Range.ListFormat.ApplyListTemplate($bulletTemplate, $true)

$item2.Range.InsertParagraphAfter()
$item3 = $item2.Next()
$item3.Style = "ListParagraph"
$item3.Range.Text = "User can manually refresh list at any time by clicking on the refresh action button on the list."
$item3.Range.ListFormat.ApplyListTemplate($bulletTemplate, $true)

# ---------------------------------------------------------------------------
# 2) Move <w:lastRenderedPageBreak/> from the run before "Toast closes and
#    nothing is added to database" to the run before "If save is successful,
#    new saved record is returned from API".
# ---------------------------------------------------------------------------

# 2a. Add the page-break marker to the "If save is successful..." paragraph.
$rngAdd = $d.Content
$rngAdd.Find.Execute("If save is successful, new saved record is returned from API", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$paraAdd = $rngAdd.Paragraphs(1)
$prngAdd = $paraAdd.Range

$xmlAdd = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="6645AF1E" w14:textId="17849057" w:rsidR="00284EB7" w:rsidRDefault="00284EB7" w:rsidP="00352563"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>If save is successful, new saved record is returned from API</w:t></w:r><w:r w:rsidR="00C614BB"><w:t>. If still on the list, the list record is replaced by the new returned record (with the database-assigned id)</w:t></w:r><w:r w:rsidR="003E2355"><w:t>, or just id is replaced</w:t></w:r><w:r w:rsidR="00C614BB"><w:t>.</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$prngAdd.InsertXML($xmlAdd)

# 2b. Remove the page-break marker from the "Toast closes..." paragraph.
$rngDel = $d.Content
$rngDel.Find.Execute("Toast closes and nothing is added to database", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$paraDel = $rngDel.Paragraphs(1)
$prngDel = $paraDel.Range

$xmlDel = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
<w:body>
<w:p w14:paraId="2C71C9F8" w14:textId="4069D832" w:rsidR="00DB4492" w:rsidRDefault="00DB4492" w:rsidP="00DB4492"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="4"/></w:numPr></w:pPr><w:r><w:t>Toast closes and nothing is added to database</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@
$prngDel.InsertXML($xmlDel)

Write-Output "Edit complete"
